$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on the Price/Volume columns so numeric-looking strings
# (e.g. "1.00", "6.88") are written back as text, matching the source data,
# then restore the original (unstyled) cell style once done.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "61.171.34"
$ws.Range("E2").Value = "  -1.93%  "
$ws.Range("D3").Value = "2.972.32"
$ws.Range("E3").Value = "  -1.62%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.27%  "
$ws.Range("D5").Value = "593.40"
$ws.Range("E5").Value = "  +1.43%  "
$ws.Range("D6").Value = "142.76"
$ws.Range("E6").Value = "  -3.26%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "0.511"
$ws.Range("E8").Value = "  -2.63%  "
$ws.Range("D9").Value = "2.967.14"
$ws.Range("E9").Value = "  -1.77%  "
$ws.Range("D10").Value = "0.147"
$ws.Range("E10").Value = "  -1.55%  "
$ws.Range("D11").Value = "6.00"
$ws.Range("E11").Value = "  +2.54%  "
$ws.Range("D12").Value = "0.450"
$ws.Range("E12").Value = "  +1.34%  "
$ws.Range("D13").Value = "0.0000225"
$ws.Range("E13").Value = "  -1.88%  "
$ws.Range("D14").Value = "33.92"
$ws.Range("E14").Value = "  -1.91%  "
$ws.Range("E15").Value = "  +3.02%  "
$ws.Range("D16").Value = "3.465.64"
$ws.Range("E16").Value = "  -1.45%  "
$ws.Range("D17").Value = "6.88"
$ws.Range("E17").Value = "  -2.68%  "
$ws.Range("D18").Value = "61.211.92"
$ws.Range("E18").Value = "  -1.65%  "
$ws.Range("D19").Value = "2.974.93"
$ws.Range("E19").Value = "  -1.19%  "
$ws.Range("D20").Value = "443.43"
$ws.Range("E20").Value = "  -4.25%  "
$ws.Range("D21").Value = "13.87"
$ws.Range("E21").Value = "  -0.71%  "
$ws.Range("D22").Value = "0.677"
$ws.Range("E22").Value = "  -1.16%  "
$ws.Range("D23").Value = "7.30"
$ws.Range("E23").Value = "  -2.06%  "
$ws.Range("D24").Value = "80.98"
$ws.Range("E24").Value = "  -0.85%  "
$ws.Range("D25").Value = "10.67"
$ws.Range("E25").Value = "  +5.72%  "
$ws.Range("D26").Value = "2.17"
$ws.Range("E26").Value = "  -3.82%  "
$ws.Range("D27").Value = "11.96"
$ws.Range("E27").Value = "  -2.83%  "
$ws.Range("E28").Value = "  +0.10%  "
$ws.Range("D29").Value = "2.69"
$ws.Range("E29").Value = "  +2.38%  "
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.38%  "
$ws.Range("D31").Value = "7.13"
$ws.Range("E31").Value = "  -0.21%  "
$ws.Range("D32").Value = "2.04"
$ws.Range("E32").Value = "  -2.60%  "
$ws.Range("D33").Value = "27.03"
$ws.Range("E33").Value = "  -5.47%  "
$ws.Range("D34").Value = "0.109"
$ws.Range("E34").Value = "  +0.11%  "
$ws.Range("D35").Value = "0.0₃0801"
$ws.Range("E35").Value = "  -0.37%  "
$ws.Range("D36").Value = "1.01"
$ws.Range("E36").Value = "  -1.55%  "
$ws.Range("D37").Value = "5.73"
$ws.Range("E37").Value = "  -0.93%  "
$ws.Range("D38").Value = "50.04"
$ws.Range("E38").Value = "  -0.67%  "
$ws.Range("D39").Value = "8.89"
$ws.Range("E39").Value = "  -2.24%  "
$ws.Range("D40").Value = "1.99"
$ws.Range("E40").Value = "  -6.00%  "
$ws.Range("E41").Value = "  +8.36%  "
$ws.Range("D42").Value = "2.80"
$ws.Range("E42").Value = "  -4.50%  "
$ws.Range("D43").Value = "384.18"
$ws.Range("E43").Value = "  -1.58%  "
$ws.Range("D44").Value = "39.15"
$ws.Range("E44").Value = "  +4.90%  "
$ws.Range("D45").Value = "0.267"
$ws.Range("E45").Value = "  -2.90%  "
$ws.Range("D46").Value = "0.0346"
$ws.Range("E46").Value = "  -2.85%  "
$ws.Range("D47").Value = "2.675.86"
$ws.Range("E47").Value = "  -2.94%  "
$ws.Range("D48").Value = "130.41"
$ws.Range("E48").Value = "  +1.33%  "
$ws.Range("D50").Value = "0.106"
$ws.Range("E50").Value = "  -2.06%  "
$ws.Range("D51").Value = "2.13"
$ws.Range("E51").Value = "  -1.38%  "

$ws.Range("D2:E51").Style = "Normal"
